# Weekly Fruta/Hortaliza update: insert 4 new price rows for
# "Vega Monumental Concepción" - Palta (Hass, Quillota) dated 2022-03-08
# (serial 44628), pushing the existing rows 463-473 down to 467-477.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 463.
$ws.Rows("463:466").Insert()

# Helper: write one full record (columns A-T) into a given row.
function Set-PaltaRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Region
    )

    $ws.Cells.Item($Row, 1).Value = 11
    $ws.Cells.Item($Row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($Row, 3).Value = "Bíobío"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = $ws.Cells.Item(462, 4).NumberFormat
    $ws.Cells.Item($Row, 5).Value = 8
    $ws.Cells.Item($Row, 6).Value = "Fruta"
    $ws.Cells.Item($Row, 7).Value = 100106
    $ws.Cells.Item($Row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($Row, 9).Value = 100106002
    $ws.Cells.Item($Row, 10).Value = "Palta"
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = "$/kilo (en caja de 15 kilos)"
    $ws.Cells.Item($Row, 18).Value = $Region
    $ws.Cells.Item($Row, 19).Value = $PrecioProm
    $ws.Cells.Item($Row, 20).Value = 1
}

Set-PaltaRow 463 44628 "Hass" "Especial" 50  3200 3200 3200 "Provincia de Quillota"
Set-PaltaRow 464 44628 "Hass" "Primera"  50  2900 2900 2900 "Provincia de Quillota"
Set-PaltaRow 465 44628 "Hass" "Segunda"  100 2600 2600 2600 "Provincia de Quillota"
Set-PaltaRow 466 44628 "Hass" "Tercera"  50  2300 2300 2300 "Provincia de Quillota"
